# Weekly update for the Cebollín / Terminal La Palmera de La Serena sheet.
# A new week's price record is inserted at row 110, which pushes every
# existing record from row 110 down through row 224 one row further down
# (the old row 224 ends up as the new row 225), growing the used range
# from A1:R224 to A1:R225.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 110; Excel shifts rows 110..224 down to 111..225
# automatically (carrying their values/styles with them), so the rest of
# the table does not need to be touched by hand.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new week's data. The
# non-varying descriptive columns (market/region/category/etc.) repeat the
# same values used throughout the rest of the sheet.
$ws.Range("A110").Value = 8
$ws.Range("B110").Value = "Terminal La Palmera de La Serena"
$ws.Range("C110").Value = "Coquimbo"
$ws.Range("D110").Value = 44781
$ws.Range("D110").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E110").Value = 4
$ws.Range("F110").Value = 100112037
$ws.Range("G110").Value = "Cebollín"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 1000
$ws.Range("K110").Value = 1400
$ws.Range("L110").Value = 1600
$ws.Range("M110").Value = 1500
$ws.Range("N110").Value = "$/paquete 6 unidades"
$ws.Range("O110").Value = "Provincia del Elquí"
$ws.Range("P110").Value = 250
$ws.Range("Q110").Value = 6
$ws.Range("R110").Value = "Hortaliza"
